$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 63, shifting existing rows 63-78 down to 64-79.
$ws.Rows.Item(63).Insert()

# Fill the newly inserted row 63 with the new record's data.
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44841
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100112026
$ws.Cells.Item(63, 7).Value = "Haba"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 30
$ws.Cells.Item(63, 11).Value = 10000
$ws.Cells.Item(63, 12).Value = 10000
$ws.Cells.Item(63, 13).Value = 10000
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 400
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
